{"js": "// Remove the four paragraphs that follow the \"LOB1208: Qu\u00edmica Anal\u00edtica\n// Ambiental I (Requisito)\" paragraph and that precede the trailing pair of\n// empty paragraphs before the section break:\n//   1) an empty paragraph\n//   2) \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n//   3) an empty paragraph\n//   4) an empty paragraph with pageBreakBefore\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst marker = \"LOB1208: Qu\u00edmica Anal\u00edtica Ambiental I (Requisito)\";\nlet markerIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === marker) {\n    markerIndex = i;\n    break;\n  }\n}\n\nif (markerIndex !== -1) {\n  // Delete the next four paragraphs (the ones immediately following the\n  // marker paragraph). Deleting from the end keeps earlier indices valid.\n  const toDelete = [];\n  for (let i = markerIndex + 1; i <= markerIndex + 4 && i < paragraphs.items.length; i++) {\n    toDelete.push(paragraphs.items[i]);\n  }\n  for (let i = toDelete.length - 1; i >= 0; i--) {\n    toDelete[i].delete();\n  }\n  await context.sync();\n}\n", "ps1": "# Remove the four paragraphs that follow the \"LOB1208: Qu\u00edmica Anal\u00edtica\n# Ambiental I (Requisito)\" paragraph and that precede the trailing pair of\n# empty paragraphs before the section break:\n#   1) an empty paragraph\n#   2) \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n#   3) an empty paragraph\n#   4) an empty paragraph with pageBreakBefore\n$d = $word.ActiveDocument\n\n$marker = \"LOB1208: Qu\u00edmica Anal\u00edtica Ambiental I (Requisito)\"\n$idx = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text.TrimEnd(\"`r\", \"`a\")\n    if ($t -eq $marker) {\n        $idx = $i\n        break\n    }\n}\n\nif ($idx -ne -1) {\n    # Delete the next four paragraphs. Walk backwards so earlier indices\n    # stay valid as each paragraph is removed.\n    $lastToDelete = [Math]::Min($idx + 4, $d.Paragraphs.Count)\n    for ($i = $lastToDelete; $i -ge $idx + 1; $i--) {\n        $d.Paragraphs.Item($i).Range.Delete()\n    }\n}\n"}
